$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commitments")

# Populate the previously-empty row 2 with the commitment record that was
# moved here as part of extracting datalock validation / commitment
# matching into its own service.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "1-001"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 42948
$ws.Range("E2").Value = 43313
$ws.Range("F2").Value = 15000
$ws.Range("G2").Value = 23
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 403
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 42948
$ws.Range("N2").Value = "NULL"
$ws.Range("O2").Value = "NULL"
$ws.Range("P2").NumberFormat = "mm:ss.0"
$ws.Range("P2").Value = 0

# Make Commitments the active/selected sheet (previously Payments was).
$ws.Activate()
$ws.Range("A2:P2").Select()
